$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2: "Meteors" - Lao column (D) and Total column (G)
$ws.Range("D2").Value = 2060
$ws.Range("G2").Value = 2060

# Row 3: "Successes" - Lao column (D) and Total column (G)
$ws.Range("D3").Value = 19
$ws.Range("G3").Value = 19
